$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows (rows 19 and 20), shifting nothing below them
$ws.Range("A19:B20").Delete()

# Update the material labels for rows 10-18 to reflect the new inventory list
$ws.Range("A10").Value = "PP"
$ws.Range("A11").Value = "PS"
$ws.Range("A12").Value = "Black chip"
$ws.Range("A13").Value = "UV Stabilizer"
$ws.Range("A14").Value = "red chip"
$ws.Range("A15").Value = "blue chip"
$ws.Range("A16").Value = "PSC"
$ws.Range("A17").Value = "green ship"
$ws.Range("A18").Value = "Slip"
